$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 00:04"

# Update country case-count rows whose data changed and/or whose
# ranking (row position, by total cases "Casos totales") shifted.
# Full A:H rows are rewritten per affected row for safety/clarity.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1342723
$ws.Cells.Item(4, 3).Value = 20938
$ws.Cells.Item(4, 4).Value = 232869
$ws.Cells.Item(4, 5).Value = 1029928
$ws.Cells.Item(4, 6).Value = 16801
$ws.Cells.Item(4, 7).Value = 1311
$ws.Cells.Item(4, 8).Value = 79926

# Row 59: Kazajistan
$ws.Cells.Item(59, 1).Value = "Kazajistan"
$ws.Cells.Item(59, 2).Value = 4975
$ws.Cells.Item(59, 3).Value = 141
$ws.Cells.Item(59, 4).Value = 1776
$ws.Cells.Item(59, 5).Value = 3168
$ws.Cells.Item(59, 6).Value = 31
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 31

# Row 115: Gabon
$ws.Cells.Item(115, 1).Value = "Gabon"
$ws.Cells.Item(115, 2).Value = 661
$ws.Cells.Item(115, 3).Value = 41
$ws.Cells.Item(115, 4).Value = 110
$ws.Cells.Item(115, 5).Value = 543
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 8

# Row 116: Kenia
$ws.Cells.Item(116, 1).Value = "Kenia"
$ws.Cells.Item(116, 2).Value = 649
$ws.Cells.Item(116, 3).Value = 28
$ws.Cells.Item(116, 4).Value = 207
$ws.Cells.Item(116, 5).Value = 412
$ws.Cells.Item(116, 6).Value = 1
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 30

# Row 117: Guinea-Bisau
$ws.Cells.Item(117, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(117, 2).Value = 641
$ws.Cells.Item(117, 3).Value = 47
$ws.Cells.Item(117, 4).Value = 25
$ws.Cells.Item(117, 5).Value = 613
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 1
$ws.Cells.Item(117, 8).Value = 3

# Row 118: San Marino
$ws.Cells.Item(118, 1).Value = "San Marino"
$ws.Cells.Item(118, 2).Value = 637
$ws.Cells.Item(118, 3).Value = 14
$ws.Cells.Item(118, 4).Value = 126
$ws.Cells.Item(118, 5).Value = 470
$ws.Cells.Item(118, 6).Value = 3
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 41

# Row 119: Georgia
$ws.Cells.Item(119, 1).Value = "Georgia"
$ws.Cells.Item(119, 2).Value = 626
$ws.Cells.Item(119, 3).Value = 3
$ws.Cells.Item(119, 4).Value = 297
$ws.Cells.Item(119, 5).Value = 319
$ws.Cells.Item(119, 6).Value = 6
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 10

# Row 192: Belice
$ws.Cells.Item(192, 1).Value = "Belice"
$ws.Cells.Item(192, 2).Value = 18
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 16
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 2

# Row 193: Nueva Caledonia
$ws.Cells.Item(193, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 18
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

# Row 205: Comoras
$ws.Cells.Item(205, 1).Value = "Comoras"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 3
$ws.Cells.Item(205, 4).Value = 0
$ws.Cells.Item(205, 5).Value = 10
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

# Row 206: Seychelles
$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 207: Montserrat
$ws.Cells.Item(207, 1).Value = "Montserrat"
$ws.Cells.Item(207, 2).Value = 11
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 8
$ws.Cells.Item(207, 5).Value = 2
$ws.Cells.Item(207, 6).Value = 1
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 1

# Row 208: Groenlandia
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 2).Value = 11
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

# Row 209: Surinam
$ws.Cells.Item(209, 1).Value = "Surinam"
$ws.Cells.Item(209, 2).Value = 10
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 9
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 1

Write-Output "applied updates"
